$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-blank GB MAKE / GB MODEL cells for rows 21-24
# (F column -> "ELECON", G column -> new "EP Series" value)
$ws.Range("F21").Value = "ELECON"
$ws.Range("G21").Value = "EP Series"

$ws.Range("F22").Value = "ELECON"
$ws.Range("G22").Value = "EP Series"

$ws.Range("F23").Value = "ELECON"
$ws.Range("G23").Value = "EP Series"

$ws.Range("F24").Value = "ELECON"
$ws.Range("G24").Value = "EP Series"

# Update the saved view: drop the frozen/scrolled topLeftCell and move the
# active selection to J20 (single cell) instead of the B2:I24 block.
[void]$ws.Range("J20").Select()
